$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.314.92"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "1.588.78"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  -0.63%  "

$ws.Range("D5").Value = "210.13"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.246"
$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0610"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").Value = "1.812.09"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "4.06"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").Value = "1.580.94"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "64.34"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "26.323.21"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("E19").Value = "  +6.00%  "

$ws.Range("D20").Value = "210.85"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("E21").Value = "  -0.75%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("E24").Value = "  -2.79%  "

$ws.Range("D25").Value = "144.45"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("E27").Value = "  -0.29%  "

$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  -0.54%  "

$ws.Range("E33").Value = "  +2.24%  "

$ws.Range("D34").Value = "1.318.98"
$ws.Range("E34").Value = "  +2.95%  "

$ws.Range("E35").Value = "  -2.18%  "

$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("E39").Value = "  -13.28%  "

$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "5.63"
$ws.Range("E42").Value = "  +4.10%  "

$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "1.724.66"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").Value = "87.48"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("E48").Value = "  -4.80%  "

$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("E50").Value = "  -4.42%  "

$ws.Range("E51").Value = "  -0.67%  "
